$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.666.12"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.589.68"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.71"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.25"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0590"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0867"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.817.25"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.592.88"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("E15").Value = "  -2.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.679.07"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.18"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.42"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.76"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.66"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("E26").Value = "  +3.62%  "
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0474"
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.374.79"
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("E38").Value = "  +2.12%  "
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.824"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.980"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.24"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("E44").Value = "  +4.21%  "
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.728.34"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("E49").Value = "  +4.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0963"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("E51").Value = "  -0.33%  "